$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1982248520710059
$ws.Range("C2").Value = 0.5502958579881657
$ws.Range("J2").Value = 0.008875739644970414
$ws.Range("P2").Value = 0.1420118343195266
$ws.Range("S2").Value = 0.1005917159763314
$ws.Range("B3").Value = 0.005319148936170213
$ws.Range("C3").Value = 0.02659574468085106
$ws.Range("J3").Value = 0.02127659574468085
$ws.Range("P3").Value = 0.8031914893617021
$ws.Range("S3").Value = 0.1436170212765958
$ws.Range("J4").Value = 0.0425531914893617
$ws.Range("P4").Value = 0.7021276595744681
$ws.Range("S4").Value = 0.2553191489361702
$ws.Range("B6").Value = 0.07981220657276995
$ws.Range("D6").Value = 0.04694835680751173
$ws.Range("F6").Value = 0.09389671361502347
$ws.Range("J6").Value = 0.1643192488262911
$ws.Range("O6").Value = 0.0187793427230047
$ws.Range("Q6").Value = 0.1877934272300469
$ws.Range("R6").Value = 0.04225352112676056
$ws.Range("S6").Value = 0.3661971830985916
$ws.Range("B7").Value = 0.1275510204081633
$ws.Range("D7").Value = 0.02551020408163265
$ws.Range("F7").Value = 0.06122448979591837
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("O7").Value = 0.02551020408163265
$ws.Range("Q7").Value = 0.2091836734693878
$ws.Range("R7").Value = 0.07142857142857142
$ws.Range("S7").Value = 0.336734693877551
$ws.Range("B8").Value = 0.09162303664921466
$ws.Range("D8").Value = 0.02094240837696335
$ws.Range("F8").Value = 0.06806282722513089
$ws.Range("J8").Value = 0.1073298429319372
$ws.Range("O8").Value = 0.02094240837696335
$ws.Range("Q8").Value = 0.2146596858638743
$ws.Range("R8").Value = 0.07591623036649214
$ws.Range("S8").Value = 0.4005235602094241
$ws.Range("B9").Value = 0.1137026239067055
$ws.Range("D9").Value = 0.01166180758017493
$ws.Range("F9").Value = 0.05830903790087463
$ws.Range("J9").Value = 0.1166180758017493
$ws.Range("O9").Value = 0.02040816326530612
$ws.Range("Q9").Value = 0.2478134110787172
$ws.Range("R9").Value = 0.05830903790087463
$ws.Range("S9").Value = 0.3731778425655977
$ws.Range("B10").Value = 0.1070911722141824
$ws.Range("D10").Value = 0.01519536903039074
$ws.Range("E10").Value = 0.0007235890014471779
$ws.Range("F10").Value = 0.05716353111432707
$ws.Range("J10").Value = 0.1143270622286541
$ws.Range("O10").Value = 0.0195369030390738
$ws.Range("Q10").Value = 0.2445730824891462
$ws.Range("R10").Value = 0.08321273516642547
$ws.Range("S10").Value = 0.3581765557163531
$ws.Range("G11").Value = 0.1151079136690648
$ws.Range("J11").Value = 0.07194244604316546
$ws.Range("K11").Value = 0.1654676258992806
$ws.Range("L11").Value = 0.6366906474820144
$ws.Range("S11").Value = 0.01079136690647482
$ws.Range("G12").Value = 0.7868852459016393
$ws.Range("J12").Value = 0.1530054644808743
$ws.Range("L12").Value = 0.03278688524590164
$ws.Range("S12").Value = 0.0273224043715847
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.06060606060606061
$ws.Range("F15").Value = 0.01716738197424893
$ws.Range("H15").Value = 0.1330472103004292
$ws.Range("I15").Value = 0.1158798283261803
$ws.Range("J15").Value = 0.3819742489270386
$ws.Range("K15").Value = 0.06437768240343347
$ws.Range("M15").Value = 0.004291845493562232
$ws.Range("O15").Value = 0.04721030042918455
$ws.Range("S15").Value = 0.2360515021459227
$ws.Range("F16").Value = 0.03056768558951965
$ws.Range("H16").Value = 0.1091703056768559
$ws.Range("I16").Value = 0.1091703056768559
$ws.Range("J16").Value = 0.4497816593886463
$ws.Range("K16").Value = 0.1266375545851528
$ws.Range("M16").Value = 0.008733624454148471
$ws.Range("N16").Value = 0.008733624454148471
$ws.Range("O16").Value = 0.06986899563318777
$ws.Range("S16").Value = 0.08733624454148471
$ws.Range("F17").Value = 0.01030927835051546
$ws.Range("H17").Value = 0.1580756013745704
$ws.Range("I17").Value = 0.1288659793814433
$ws.Range("J17").Value = 0.4329896907216495
$ws.Range("K17").Value = 0.08075601374570447
$ws.Range("M17").Value = 0.02920962199312715
$ws.Range("N17").Value = 0.001718213058419244
$ws.Range("O17").Value = 0.06185567010309279
$ws.Range("S17").Value = 0.09621993127147767
$ws.Range("F18").Value = 0.005376344086021506
$ws.Range("H18").Value = 0.1451612903225807
$ws.Range("I18").Value = 0.1612903225806452
$ws.Range("J18").Value = 0.4623655913978494
$ws.Range("K18").Value = 0.08602150537634409
$ws.Range("M18").Value = 0.01075268817204301
$ws.Range("O18").Value = 0.04301075268817205
$ws.Range("S18").Value = 0.08602150537634409
$ws.Range("F19").Value = 0.01108471892319873
$ws.Range("H19").Value = 0.1670625494853523
$ws.Range("I19").Value = 0.1456848772763262
$ws.Range("J19").Value = 0.3966745843230404
$ws.Range("K19").Value = 0.09501187648456057
$ws.Range("M19").Value = 0.009501187648456057
$ws.Range("O19").Value = 0.06730007917656373
$ws.Range("S19").Value = 0.107680126682502
